# "Better labels on Mapping entities"
# Appends a running row index (" - N") to the generated mapping label in
# column N of the SSSOM sheet, and leaves the sheet/selection state the way
# the author left it after making the edit (SSSOM tab active, N2 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SSSOM")

$firstRow = 2
$lastRow = 189

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 14)
    $cell.Formula = '=CONCATENATE(B' + $r + ', " - mapping to IUCN GET - ", ROW(B' + $r + ')-1)'
}

# Author ended the edit with the SSSOM sheet active (was "header" before)
# and cell N2 selected, scrolled so column I is the leftmost visible column.
$ws.Activate()
$ws.Range("I1").Select()
$ws.Range("N2").Select()
